# Add new "Lung Cancer Dataset" row (row 11) to Sheet1, mirroring the
# formatting pattern used by the existing dataset rows, and refresh a
# couple of cosmetic view settings that changed alongside it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new hyperlink first (on the still-empty B11) so Excel's automatic
# "Hyperlink" cell style doesn't clobber formatting we paste afterwards ---
$lungUrl = "https://www.kaggle.com/datasets/thedevastator/cancer-patients-and-air-pollution-a-new-link/"
$ws.Hyperlinks.Add($ws.Range("B11"), $lungUrl) | Out-Null

# --- clone row 10's cell formatting onto row 11 (same visual pattern:
# name / link / detection / features columns) ---
$ws.Range("A10:D10").Copy() | Out-Null
$ws.Range("A11:D11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- values for the new dataset row ---
$ws.Range("A11").Value = "Lung Cancer Dataset"
$ws.Range("C11").Value = "Severiaity of Lung Cancer"
$ws.Range("D11").Value = "Age-	The age of the patient. (Numeric)
Gender-	The gender of the patient. (Categorical)
Air Pollution-	The level of air pollution exposure of the patient. (Categorical)
Alcohol use-	The level of alcohol use of the patient. (Categorical)
Dust Allergy-	The level of dust allergy of the patient. (Categorical)
Occupational Hazards-	The level of occupational hazards of the patient. (Categorical)
Genetic Risk-	The level of genetic risk of the patient. (Categorical)
chronic Lung Disease-	The level of chronic lung disease of the patient. (Categorical)
Balanced Diet-	The level of balanced diet of the patient. (Categorical)
Obesity-	The level of obesity of the patient. (Categorical)
Smoking-	The level of smoking of the patient. (Categorical)
Passive Smoker-	The level of passive smoker of the patient. (Categorical)
Chest Pain-	The level of chest pain of the patient. (Categorical)
Coughing of Blood-	The level of coughing of blood of the patient. (Categorical)
Fatigue-	The level of fatigue of the patient. (Categorical)
Weight Loss-	The level of weight loss of the patient. (Categorical)
Shortness of Breath-	The level of shortness of breath of the patient. (Categorical)
Wheezing-	The level of wheezing of the patient. (Categorical)
Swallowing Difficulty-	The level of swallowing difficulty of the patient. (Categorical)
Clubbing of Finger Nails-	The level of clubbing of finger nails of the patient. (Categorical)
"

# --- widen name/link columns to fit the new content ---
$ws.Columns.Item(2).ColumnWidth = 20.16666666666667
$ws.Columns.Item(3).ColumnWidth = 28.799479166666668

# --- row height for the tall new row (matches the other wrapped-text rows) ---
$ws.Rows.Item(11).RowHeight = 319

# --- view tidy-up: drop the stale scroll/active-cell position and select
# the title bar merged cell instead ---
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A1:E1").Select()
